# Adding stratification violation to the data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the data values that introduce the stratification violation ---
$ws.Range("B56").Value = -2.2
$ws.Range("B57").Value = 2.02737541784593
$ws.Range("B61").Value = 2.2
$ws.Range("B62").Value = -1.2
$ws.Range("B65").Value = 2.2
$ws.Range("B67").Value = -2.2

# --- Update the sheet's default column width (best effort) ---
$ws.StandardWidth = 11.66015625

# --- Update the view / scroll position / selection to match the saved state ---
$win = $excel.ActiveWindow
$win.ScrollColumn = 1
$win.ScrollRow = 37
$ws.Range("B68").Select()
